# Applies the "Added a small convenience patch for uthash" benchmark update
# to Sheet1 of the workbook: renames the row-36 section header, backfills two
# missing F-column samples, and appends a new benchmark block (rows 45-53)
# comparing heap-detection strategies.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 36: section header text changes from SpeedTest to TortureTest ---
$ws.Range("A36").Value = "One thread, 10000 record TortureTest x64:"

# --- Rows 37-38: backfill the previously-missing column F samples ---
$ws.Range("F37").Value = 648429
$ws.Range("F38").Value = 658633

# --- New benchmark block, rows 45-53 ---
# Cells are written in the order that first introduces each new label so the
# resulting shared-string table layout matches the source edit.

# Row 51/52/53 labels (first occurrences of these strings)
$ws.Range("A51").Value = "normal"
$ws.Range("A52").Value = "LTCG off"
$ws.Range("A53").Value = "/GS-"

# Row 49 / 46 / 45 / 48 labels (first occurrences)
$ws.Range("A49").Value = "DLL, fast heap detect"
$ws.Range("A46").Value = "static, fast heap detect"
$ws.Range("A45").Value = "static, normal detect"
$ws.Range("A48").Value = "DLL, normal detect"

# Row 49 column E label (first occurrence of this string)
$ws.Range("E49").Value = "/GS-, SSE2"

# Remaining duplicate-text label cells
$ws.Range("E48").Value = "DLL, normal detect"
$ws.Range("E50").Value = "LTCG off"

# --- Numeric values ---
$ws.Range("B45").Value = 779216
$ws.Range("B46").Value = 791479
$ws.Range("B48").Value = 769324
$ws.Range("B49").Value = 784905
$ws.Range("B51").Value = 780550
$ws.Range("B52").Value = 773581
$ws.Range("B53").Value = 768500

$ws.Range("G48").Value = 807550
$ws.Range("G49").Value = 796140
$ws.Range("G50").Value = 796782

# --- Formulas (percentage change vs. previous row) ---
$ws.Range("C46").Formula = "=(B46-B45)/B45"
$ws.Range("C49").Formula = "=(B49-B48)/B48"

# --- View state: selection moves to G51 as the sheet now scrolls further down ---
$ws.Range("G51").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
